# "re-run RU 1001; without crop"
# The Saudi Arabia column (C) now produces #NUM! errors for every country
# row (the re-run without the crop/clip step pushes its inputs out of the
# valid domain), and a handful of other cells shift by a single ULP because
# the underlying random draw was re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Saudi Arabia") -> #NUM! for every data row (rows 2-14)
$ws.Range("C2").Value  = "#NUM!"
$ws.Range("C3").Value  = "#NUM!"
$ws.Range("C4").Value  = "#NUM!"
$ws.Range("C5").Value  = "#NUM!"
$ws.Range("C6").Value  = "#NUM!"
$ws.Range("C7").Value  = "#NUM!"
$ws.Range("C8").Value  = "#NUM!"
$ws.Range("C9").Value  = "#NUM!"
$ws.Range("C10").Value = "#NUM!"
$ws.Range("C11").Value = "#NUM!"
$ws.Range("C12").Value = "#NUM!"
$ws.Range("C13").Value = "#NUM!"
$ws.Range("C14").Value = "#NUM!"

# Tiny (last-digit) re-run precision drifts in unrelated cells
$ws.Range("I2").Value  = 0.0778600772704492
$ws.Range("L3").Value  = 0.0718945015844923
$ws.Range("L4").Value  = 0.0632360925971079
$ws.Range("G6").Value  = 0.095665580015614
$ws.Range("E9").Value  = 0.0628958989764233
$ws.Range("I11").Value = 0.0863194602307682
$ws.Range("K12").Value = 0.0706362559659104
$ws.Range("B13").Value = 0.0853428169405929
$ws.Range("F13").Value = 0.0984951828658039
$ws.Range("J13").Value = 0.0912410141275855
$ws.Range("E14").Value = 0.0815049427224147
$ws.Range("G14").Value = 0.0663133141155384
